$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 351 (shifts existing rows 351:375 down to 352:376)
$ws.Rows.Item(351).Insert()

# Populate the newly inserted row with the new weekly price record
$ws.Range("A351").Value = 7
$ws.Range("B351").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C351").Value = "Ñuble"
$ws.Range("D351").Value = 44585
$ws.Range("E351").Value = 16
$ws.Range("F351").Value = 100112020
$ws.Range("G351").Value = "Tomate"
$ws.Range("H351").Value = "Larga vida"
$ws.Range("I351").Value = "Primera"
$ws.Range("J351").Value = 600
$ws.Range("K351").Value = 6000
$ws.Range("L351").Value = 6500
$ws.Range("M351").Value = 6250
$ws.Range("N351").Value = "`$/caja 15 kilos"
$ws.Range("O351").Value = "Región del Maule"
$ws.Range("P351").Value = 417
$ws.Range("Q351").Value = 15
$ws.Range("R351").Value = "Hortaliza"
